# Applies updated crypto price (D) and 1h volume change (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.118.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.117.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.113.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.632.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.133.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.117.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.845.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "383.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "137.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("E51").Value = "  -0.30%  "
